$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 1038, shifting existing row 1038
# (and everything below it) down by one row.
$ws.Rows.Item(1038).Insert()

# Populate the newly inserted row 1038 with the new data point.
$ws.Cells.Item(1038, 1).Value = 3
$ws.Cells.Item(1038, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1038, 3).Value = "Coquimbo"
$ws.Cells.Item(1038, 4).Value = 45132
$ws.Cells.Item(1038, 5).Value = 5
$ws.Cells.Item(1038, 6).Value = 100114001
$ws.Cells.Item(1038, 7).Value = "Papa"
$ws.Cells.Item(1038, 8).Value = "Rosara"
$ws.Cells.Item(1038, 9).Value = "1a (guarda)"
$ws.Cells.Item(1038, 10).Value = 370
$ws.Cells.Item(1038, 11).Value = 17000
$ws.Cells.Item(1038, 12).Value = 18000
$ws.Cells.Item(1038, 13).Value = 17568
$ws.Cells.Item(1038, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(1038, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(1038, 16).Value = 703
$ws.Cells.Item(1038, 17).Value = 25
$ws.Cells.Item(1038, 18).Value = "Hortaliza"
